$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the retired catalogue entry (id JUdOAAAACAAJ - "The Lord of the Rings Sketchbook").
# Deleting row 9 shifts rows 10-31 up by one (new rows 9-30).
$ws.Rows(9).Delete()

# The row that shifted from 14 -> 13 had a blank "authors" cell; make sure it stays blank
# (ClearContents removes the cell entirely, matching the empty inlineStr in the source).
$ws.Range("E13").ClearContents()

# The row that shifted from 17 -> 16 also had a blank "authors" cell.
$ws.Range("E16").ClearContents()

# Add the newly catalogued entry (id eqPUjwEACAAJ) in its sorted position at row 22,
# pushing the remaining rows (formerly 22-31) down to 23-31.
$ws.Rows(22).Insert()

# Force text formatting on the new row's date cell so Excel keeps the published date as
# plain text instead of converting it to a date serial number.
$ws.Range("D22").NumberFormat = "@"

$ws.Range("A22").Value = 'eqPUjwEACAAJ'
$ws.Range("B22").Value = 'El Senor de Los Anillos'
$ws.Range("C22").Value = 'Titulo: El Señor de los anillos: La comunidad del anillo.Autor: J. R. R. Tolkien.Fecha de publicación: 1954.En la adormecida e idílica Comarca, un joven hobbit recibe un encargo: custodiar el Anillo Único y emprender el viaje para su destrucción en las Grietas del Destino. Acompañado por magos, hombres, elfos y enanos, atravesará la Tierra Media y se internará en las sombras de Mordor, perseguido siempre por las huestes de Sauron, el Señor Oscuro, dispuesto a recuperar su creación para establecer el dominio definitivo del Mal.'
$ws.Range("D22").Value = '2016-05-31'
$ws.Range("E22").Value = 'J. R. R. Tolkien'
